# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to match the freshly re-scraped source data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- Sheet "展览" (sheet1) : column F updates ----
$exhibitUpdates = @{
    3  = 1369
    4  = 0
    5  = 763
    9  = 0
    10 = 0
    11 = 711
    12 = 230
    13 = 25
    14 = 139
    15 = 88
    16 = 228
    19 = 334
    20 = 0
    23 = 44
    24 = 50
    25 = 0
    26 = 1054
    27 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 36
    33 = 0
    36 = 12390
    37 = 1311
    38 = 44
    40 = 51
    42 = 323
    44 = 0
    45 = 95
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# ---- Sheet "全部类型" (sheet4) : column F updates ----
$allTypeUpdates = @{
    2  = 194
    3  = 0
    4  = 19261
    5  = 763
    6  = 0
    8  = 3
    9  = 0
    10 = 469
    13 = 0
    14 = 139
    15 = 88
    16 = 0
    17 = 178
    18 = 0
    19 = 0
    22 = 42
    24 = 50
    25 = 0
    28 = 1
    29 = 0
    31 = 547
    32 = 0
    34 = 32
    38 = 0
    39 = 0
    40 = 44
    41 = 2
    43 = 247
    44 = 0
    46 = 315
    47 = 0
}

foreach ($row in $allTypeUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allTypeUpdates[$row]
}
